$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Rename "Table S2 - PERMANOVA all" -> "Table S2 - PERMANOVA"
# ------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Table S2 - PERMANOVA all")
$ws2.Name = "Table S2 - PERMANOVA"

# ------------------------------------------------------------------
# 2. Fix up a handful of P-values on Table S2
# ------------------------------------------------------------------
$ws2.Range("F3").Value = 0.00533
$ws2.Range("F4").Value = 0.00333
$ws2.Range("F8").Value = 0.002
$ws2.Range("F12").Value = 0.81945

# ------------------------------------------------------------------
# 3. Add a new "Table S5 - Species PERMANOVA" sheet after Table S4
# ------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("Table S4 - Plasticity ANOVA")
$ws5 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws4)
$ws5.Name = "Table S5 - Species PERMANOVA"

# Column widths matching the authored worksheet
$ws5.Columns.Item(1).ColumnWidth = 24.71
$ws5.Columns.Item(2).ColumnWidth = 3.71
$ws5.Columns.Item(3).ColumnWidth = 14.71
$ws5.Columns.Item(4).ColumnWidth = 5.71
$ws5.Columns.Item(5).ColumnWidth = 5.71
$ws5.Columns.Item(6).ColumnWidth = 7.71

# Header row
$ws5.Cells.Item(1, 1).Value = " "
$ws5.Cells.Item(1, 2).Value = "Df"
$ws5.Cells.Item(1, 3).Value = "Sum of Squares"
$ws5.Cells.Item(1, 4).Value = "R2"
$ws5.Cells.Item(1, 5).Value = "F"
$ws5.Cells.Item(1, 6).Value = "P-value"

# Row 2 - pCO2
$ws5.Cells.Item(2, 1).Value = "pCO2"
$ws5.Cells.Item(2, 2).Value = 3
$ws5.Cells.Item(2, 3).Value = 441584
$ws5.Cells.Item(2, 4).Value = 0.084
$ws5.Cells.Item(2, 5).Value = 14.9
$ws5.Cells.Item(2, 6).Value = 0.00067

# Row 3 - temperature
$ws5.Cells.Item(3, 1).Value = "temperature"
$ws5.Cells.Item(3, 2).Value = 1
$ws5.Cells.Item(3, 3).Value = 44470
$ws5.Cells.Item(3, 4).Value = 0.008
$ws5.Cells.Item(3, 5).Value = 4.5
$ws5.Cells.Item(3, 6).Value = 0.02398

# Row 4 - reef environment
$ws5.Cells.Item(4, 1).Value = "reef environment"
$ws5.Cells.Item(4, 2).Value = 1
$ws5.Cells.Item(4, 3).Value = 69064
$ws5.Cells.Item(4, 4).Value = 0.013
$ws5.Cells.Item(4, 5).Value = 6.99
$ws5.Cells.Item(4, 6).Value = 0.006

# Row 5 - species
$ws5.Cells.Item(5, 1).Value = "species"
$ws5.Cells.Item(5, 2).Value = 2
$ws5.Cells.Item(5, 3).Value = 1690024
$ws5.Cells.Item(5, 4).Value = 0.32
$ws5.Cells.Item(5, 5).Value = 85.53
$ws5.Cells.Item(5, 6).Value = 0.00067

# Row 6 - temperature:species
$ws5.Cells.Item(6, 1).Value = "temperature:species"
$ws5.Cells.Item(6, 2).Value = 2
$ws5.Cells.Item(6, 3).Value = 657172
$ws5.Cells.Item(6, 4).Value = 0.124
$ws5.Cells.Item(6, 5).Value = 33.26
$ws5.Cells.Item(6, 6).Value = 0.00067

# Row 7 - pCO2:species
$ws5.Cells.Item(7, 1).Value = "pCO2:species"
$ws5.Cells.Item(7, 2).Value = 6
$ws5.Cells.Item(7, 3).Value = 129807
$ws5.Cells.Item(7, 4).Value = 0.025
$ws5.Cells.Item(7, 5).Value = 2.19
$ws5.Cells.Item(7, 6).Value = 0.01999

# Row 8 - reef environment:species
$ws5.Cells.Item(8, 1).Value = "reef environment:species"
$ws5.Cells.Item(8, 2).Value = 2
$ws5.Cells.Item(8, 3).Value = 134579
$ws5.Cells.Item(8, 4).Value = 0.025
$ws5.Cells.Item(8, 5).Value = 6.81
$ws5.Cells.Item(8, 6).Value = 0.00067

# Row 9 - Residual (E/F intentionally left blank)
$ws5.Cells.Item(9, 1).Value = "Residual"
$ws5.Cells.Item(9, 2).Value = 214
$ws5.Cells.Item(9, 3).Value = 2114166
$ws5.Cells.Item(9, 4).Value = 0.4

# Row 10 - Total (E/F intentionally left blank)
$ws5.Cells.Item(10, 1).Value = "Total"
$ws5.Cells.Item(10, 2).Value = 231
$ws5.Cells.Item(10, 3).Value = 5280867
$ws5.Cells.Item(10, 4).Value = 1
